# Rename the inline picture shapes (logos) in the document's header/footer
# stories so the wp:docPr / underlying picture "name" metadata matches the
# target revision:
#   - footer (first page)  "PearsonLogo" image: image2.png -> image1.png
#   - footer (default)     "PearsonLogo" image: image2.png -> image1.png
#   - header (first page)  "BTec_Logo-Orange" image: image1.jpg -> image2.jpg

$d = $word.ActiveDocument
$sec = $d.Sections.Item(1)

# wdHeaderFooterPrimary = 1 (default header/footer)
# wdHeaderFooterFirstPage = 2 (first-page header/footer, used here because
# the section has titlePg set)

$footerDefault = $sec.Footers.Item(1)
$footerFirst   = $sec.Footers.Item(2)
$headerFirst   = $sec.Headers.Item(2)

# Pearson logo in the default (non-first-page) footer.
if ($footerDefault.Range.InlineShapes.Count -ge 1) {
    $shape = $footerDefault.Range.InlineShapes.Item(1)
    $shape.Name = "image1.png"
}

# Pearson logo in the first-page footer.
if ($footerFirst.Range.InlineShapes.Count -ge 1) {
    $shape = $footerFirst.Range.InlineShapes.Item(1)
    $shape.Name = "image1.png"
}

# BTEC logo in the first-page header.
if ($headerFirst.Range.InlineShapes.Count -ge 1) {
    $shape = $headerFirst.Range.InlineShapes.Item(1)
    $shape.Name = "image2.jpg"
}
